# Update crypto price (D) and 1h volume-change (E) figures to the latest
# scrape, per the GitHub Actions refresh job.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'67.288.05"
$ws.Range("E2").Value = "  +1.27%  "
$ws.Range("D3").Value = "'2.529.45"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("E4").Value = "  -0.01%  "
$ws.Range("D5").Value = "'591.50"
$ws.Range("E5").Value = "  +1.47%  "
$ws.Range("D6").Value = "'173.96"
$ws.Range("E6").Value = "  +4.62%  "
$ws.Range("E7").Value = "  +0.03%  "
$ws.Range("E8").Value = "  +0.28%  "
$ws.Range("D9").Value = "'2.528.51"
$ws.Range("E9").Value = "  -2.00%  "
$ws.Range("E10").Value = "  +1.48%  "
$ws.Range("E11").Value = "  +2.22%  "
$ws.Range("E12").Value = "  -0.03%  "
$ws.Range("E13").Value = "  -3.34%  "
$ws.Range("D14").Value = "'26.59"
$ws.Range("E14").Value = "  -0.17%  "
$ws.Range("D15").Value = "'2.993.43"
$ws.Range("E15").Value = "  -1.80%  "
$ws.Range("D16").Value = "'0.0000177"
$ws.Range("E16").Value = "  -0.03%  "
$ws.Range("D17").Value = "'67.083.89"
$ws.Range("E17").Value = "  +1.06%  "
$ws.Range("D18").Value = "'2.504.71"
$ws.Range("E18").Value = "  -3.25%  "
$ws.Range("E19").Value = "  +4.96%  "
$ws.Range("E20").Value = "  +0.02%  "
$ws.Range("D21").Value = "'355.36"
$ws.Range("E21").Value = "  +1.48%  "
$ws.Range("D22").Value = "'4.19"
$ws.Range("E22").Value = "  -0.70%  "
$ws.Range("E23").Value = "  +0.66%  "
$ws.Range("E24").Value = "  +6.51%  "
$ws.Range("E25").Value = "  +0.04%  "
$ws.Range("D26").Value = "'69.90"
$ws.Range("E26").Value = "  +1.65%  "
$ws.Range("D27").Value = "'9.93"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("D28").Value = "'0.999"
$ws.Range("E28").Value = "  -0.11%  "
$ws.Range("D29").Value = "'2.651.03"
$ws.Range("E29").Value = "  -2.30%  "
$ws.Range("D30").Value = "'0.0₃0982"
$ws.Range("E30").Value = "  -0.03%  "
$ws.Range("D31").Value = "'539.93"
$ws.Range("E31").Value = "  +1.89%  "
$ws.Range("D32").Value = "'8.21"
$ws.Range("E32").Value = "  +2.01%  "
$ws.Range("E33").Value = "  +1.05%  "
$ws.Range("E34").Value = "  +0.48%  "
$ws.Range("E35").Value = "  -0.87%  "
$ws.Range("E36").Value = "  +0.01%  "
$ws.Range("E37").Value = "  +0.55%  "
$ws.Range("D38").Value = "'156.38"
$ws.Range("E38").Value = "  -0.06%  "
$ws.Range("D39").Value = "'18.68"
$ws.Range("E39").Value = "  -0.17%  "
$ws.Range("D40").Value = "'18.45"
$ws.Range("E40").Value = "  +0.94%  "
$ws.Range("E41").Value = "  -1.63%  "
$ws.Range("E42").Value = "  +1.64%  "
$ws.Range("E43").Value = "  +1.03%  "
$ws.Range("E44").Value = "  +0.03%  "
$ws.Range("E45").Value = "  +3.98%  "
$ws.Range("E46").Value = "  +0.43%  "
$ws.Range("E47").Value = "  -1.17%  "
$ws.Range("E48").Value = "  -2.53%  "
$ws.Range("E49").Value = "  -0.26%  "
$ws.Range("D50").Value = "'1.70"
$ws.Range("E50").Value = "  -0.21%  "
$ws.Range("E51").Value = "  -0.14%  "
